# GLB 06-26 Hypatia Update
# Row 171 ('Dream Interpretation' / Hecate dispatch) is replaced by a new
# 'Academic Exchange' / Hypatia dispatch, and three rows are inserted after it:
#   172 - Tomb Raiding Professional / Lysandra
#   173 - DisCity Appraisal / Lysandra + McQueen
#   174 - the original Dream Interpretation / Hecate dispatch, shifted down

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Dispatch")

function Apply-RowData {
    param($ws, $rowNum, $rowData)
    foreach ($col in $rowData.Keys) {
        $val = $rowData[$col]
        $cellRef = "$col$rowNum"
        if ($val -is [string] -and $val -match '^-?[0-9]+(\.[0-9]+)?$') {
            # Force number-looking text (e.g. reward quantity "1.0") to stay text
            $ws.Range($cellRef).Value = "'" + $val
        } else {
            $ws.Range($cellRef).Value = $val
        }
    }
}

# Insert three fresh rows right after row 171, pushing nothing else around
# since row 171 is currently the last row in the sheet.
$ws.Range("A172:A174").EntireRow.Insert()

$row171 = [ordered]@{
    "A" = 1053001
    "B" = "Purple"
    "C" = "Academic Exchange"
    "D" = "学術交流"
    "E" = "학술 교류"
    "F" = "学术交流"
    "G" = "學術交流"
    "H" = "The Bureau is planning to organize an internal scientific conference to share cutting-edge developments from the scientific community, and has extended a special invitation to a particular Sinner."
    "I" = "最近、学界の第一線の情報を共有するため、管理局は内部科学交流会の開催を計画しており、特別にとあるコンビクトを招待した。"
    "J" = "최근 학계의 최신 정보를 얻기 위해 관리국은 내부 과학 교류회를 열기로 했고, 특별히 한 수감자를 초청했다."
    "K" = "近日，管理局计划组织一次内部科学交流会，分享学界的前沿资讯，特别邀请了某位禁闭者参加。"
    "L" = "近日，管理局計畫組織一次內部科學交流會，分享學界的最新資訊，特別邀請了某位禁閉者參加。"
    "M" = "Hypatia"
    "N" = "ヒパティア"
    "O" = "히파티아"
    "P" = "伊帕希娅"
    "Q" = "伊帕希婭"
    "W" = "Arsenopyrite Crystal"
    "X" = "毒砂結晶"
    "Y" = "독모래 결정"
    "Z" = "毒砂晶"
    "AA" = "毒砂晶"
    "AB" = "1.0"
    "AC" = "Arsenopyrite Concentrate"
    "AD" = "毒砂の精鉱"
    "AE" = "정교한 독모래 광석"
    "AF" = "毒砂精矿"
    "AG" = "毒砂精礦"
    "AH" = "1.0"
}

$row172 = [ordered]@{
    "A" = 1053002
    "B" = "Green"
    "C" = "Tomb Raiding Professional"
    "D" = "盗掘職人"
    "E" = "도굴 장인"
    "F" = "倒斗匠人"
    "G" = "盜墓匠人"
    "H" = "While patrolling the outskirts of DisCity, the patrol discovered a suspicious pit. The Bureau now needs to send a Sinner with relevant excavation experience to assist with the investigation."
    "I" = "パトロール隊がディスシティの外周を巡回中、不審な深い穴を発見した。管理局は現在、調査に協力するため、関連作業の経験を持つコンビクトを派遣する必要がある。"
    "J" = "순찰대가 디스시티 외곽을 순찰하던 중, 의문스러운 구덩이 하나를 발견했다. 관리국에서 관련 작업 경험이 있는 수감자를 파견해 조사 작업을 도와야 한다."
    "K" = "巡逻队于狄斯城外围巡查时发现一处可疑深坑，现需管理局派出一名有相关作业经验的禁闭者协助探查工作。"
    "L" = "巡邏隊於狄斯城外圍巡查時發現一處可疑深坑，現需管理局派出一名有相關作業經驗的禁閉者協助探查工作。"
    "M" = "Lysandra"
    "N" = "リサロ"
    "O" = "리산드라"
    "P" = "黎莎洛"
    "Q" = "黎莎洛"
    "W" = "Arsenopyrite Concentrate"
    "X" = "毒砂の精鉱"
    "Y" = "정교한 독모래 광석"
    "Z" = "毒砂精矿"
    "AA" = "毒砂精礦"
    "AB" = "1.0"
    "AC" = "Arsenopyrite Raw Ore"
    "AD" = "毒砂の原鉱"
    "AE" = "거친 독모래 광석"
    "AF" = "毒砂粗矿"
    "AG" = "毒砂粗礦"
    "AH" = "1.0"
}

$row173 = [ordered]@{
    "A" = 1053003
    "B" = "Blue"
    "C" = "DisCity Appraisal"
    "D" = "ディス宝物鑑定"
    "E" = "디스 보물 감정"
    "F" = "狄斯鉴宝"
    "G" = "狄斯鑒寶"
    "H" = "The Public Security Bureau recently confiscated a collection of ancient artifacts contaminated with Mania from the black market. Some fakes are mixed among them, which must be destroyed immediately. They require assistance of several Sinners skilled in artifact and artwork appraisal."
    "I" = "最近、治安局は闇市で狂瞳病に汚染された古代美術品を押収した。その中には贋作も混じっており、直接破棄する必要がある。美術品や古代遺物の鑑定に長けたコンビクト数名の協力が必要だ。"
    "J" = "치안국은 최근 암시장에서 변이에 오염된 고대 예술품을 압수했고, 그중에는 위조품도 섞여 있어 즉시 폐기 조치가 필요하다. 이에 따라 유물 및 예술품 감정에 능숙한 수감자 몇 명의 협조가 필요하다."
    "K" = "治安局近日在黑市缴获一批受狂厄污染的古代艺术品，其中混有赝品，需直接销毁，现需几名擅长文物及艺术品鉴定的禁闭者协助工作。"
    "L" = "治安局近日在黑市繳獲一批受狂厄汙染的古代藝術品，其中混有贗品，需直接銷毀，現需幾名擅長文物及藝術品鑑定的禁閉者協助工作。"
    "M" = "Lysandra"
    "N" = "リサロ"
    "O" = "리산드라"
    "P" = "黎莎洛"
    "Q" = "黎莎洛"
    "R" = "McQueen"
    "S" = "マックイーン"
    "T" = "맥퀸"
    "U" = "麦昆"
    "V" = "麥昆"
    "W" = "Arsenopyrite Concentrate"
    "X" = "毒砂の精鉱"
    "Y" = "정교한 독모래 광석"
    "Z" = "毒砂精矿"
    "AA" = "毒砂精礦"
    "AB" = "1.0"
    "AC" = "Arsenopyrite Concentrate"
    "AD" = "毒砂の精鉱"
    "AE" = "정교한 독모래 광석"
    "AF" = "毒砂精矿"
    "AG" = "毒砂精礦"
    "AH" = "1.0"
}

$row174 = [ordered]@{
    "A" = 1099001
    "B" = "Green"
    "C" = "Dream Interpretation"
    "D" = "夢の解析"
    "E" = "꿈의 해석"
    "F" = "梦的解析"
    "G" = "夢的解析"
    "H" = "Residents in some areas of Eastside often suffer from nightmares. It is necessary to find out why."
    "I" = "ニューシティの一部エリアの住民がよくナイトメアにうなされている。具体的な原因を調査しなければならない。"
    "J" = "신성 일부 지역 주민들은 악몽을 자주꾼다. 구체적인 원인을 철저히 조사해야 한다."
    "K" = "新城一些区域的居民经常做噩梦，需要查清具体原因。"
    "L" = "新城一些區域的居民經常做惡夢，需要查明具體原因。"
    "M" = "Hecate"
    "N" = "ヘカテー"
    "O" = "헤카테"
    "P" = "赫卡蒂"
    "Q" = "赫卡蒂"
    "W" = "Arsenopyrite Concentrate"
    "X" = "毒砂の精鉱"
    "Y" = "정교한 독모래 광석"
    "Z" = "毒砂精矿"
    "AA" = "毒砂精礦"
    "AB" = "1.0"
    "AC" = "Arsenopyrite Raw Ore"
    "AD" = "毒砂の原鉱"
    "AE" = "거친 독모래 광석"
    "AF" = "毒砂粗矿"
    "AG" = "毒砂粗礦"
    "AH" = "1.0"
}

Apply-RowData $ws 171 $row171
Apply-RowData $ws 172 $row172
Apply-RowData $ws 173 $row173
Apply-RowData $ws 174 $row174

"Row 171 A/C/M: $($ws.Range('A171').Value()) / $($ws.Range('C171').Value()) / $($ws.Range('M171').Value())"
"Row 172 A/C/M: $($ws.Range('A172').Value()) / $($ws.Range('C172').Value()) / $($ws.Range('M172').Value())"
"Row 173 A/C/M: $($ws.Range('A173').Value()) / $($ws.Range('C173').Value()) / $($ws.Range('M173').Value())"
"Row 174 A/C/M: $($ws.Range('A174').Value()) / $($ws.Range('C174').Value()) / $($ws.Range('M174').Value())"
